$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.799.36'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.129.56'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '532.27'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').Value = '138.60'
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.127.12'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('E9').Value = '  +4.73%  '
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  +4.80%  '
$ws.Range('D13').Value = '3.669.87'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').Value = '25.55'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').Value = '57.958.57'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '3.127.79'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '6.03'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').Value = '12.71'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '8.12'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').Value = '359.77'
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '68.93'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('D29').Value = '7.30'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').Value = '1.87'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = '6.09'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').Value = '21.40'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').Value = '5.11'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('D34').Value = '1.14'
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').Value = '158.46'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').Value = '6.07'
$ws.Range('E36').Value = '  -1.15%  '
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('D39').Value = '1.66'
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').Value = '2.499.76'
$ws.Range('E41').Value = '  +6.81%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '0.698'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '4.00'
$ws.Range('E43').Value = '  -4.25%  '
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('D45').Value = '3.175.32'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '0.0268'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').Value = '0.988'
$ws.Range('E48').Value = '  +3.13%  '
$ws.Range('D49').Value = '6.07'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = '19.77'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = '0.741'
$ws.Range('E51').Value = '  -2.80%  '
